$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("D3").Value = "08/20/2021"

# Row 4
$ws.Range("D4").Value = "09/09/2021"
$ws.Range("M4").Value = 100

# Row 5
$ws.Range("D5").Value = "08/02/2021"
$ws.Range("M5").Value = 200
$ws.Range("N5").Value = 20000
$ws.Range("O5").Value = 21000
$ws.Range("P5").Value = 20500
$ws.Range("S5").Value = 1025

# Row 6
$ws.Range("D6").Value = "06/17/2021"
$ws.Range("M6").Value = 140

# Row 7
$ws.Range("D7").Value = "08/19/2021"
$ws.Range("M7").Value = 200
$ws.Range("N7").Value = 20000
$ws.Range("O7").Value = 21000
$ws.Range("P7").Value = 20500
$ws.Range("S7").Value = 1025

# Row 8
$ws.Range("D8").Value = "08/26/2021"
$ws.Range("M8").Value = 100
$ws.Range("N8").Value = 20000
$ws.Range("O8").Value = 21000
$ws.Range("P8").Value = 20500
$ws.Range("S8").Value = 1025

# Row 9
$ws.Range("D9").Value = "08/09/2021"
$ws.Range("M9").Value = 160

# Row 10
$ws.Range("D10").Value = "09/02/2021"
$ws.Range("M10").Value = 160
$ws.Range("N10").Value = 20000
$ws.Range("O10").Value = 21000
$ws.Range("P10").Value = 20500
$ws.Range("S10").Value = 1025

# Row 11
$ws.Range("D11").Value = "08/27/2021"
$ws.Range("M11").Value = 260
$ws.Range("O11").Value = 22000
$ws.Range("P11").Value = 21115
$ws.Range("S11").Value = 1056

# Row 13
$ws.Range("D13").Value = "09/03/2021"
$ws.Range("M13").Value = 140
$ws.Range("N13").Value = 20000
$ws.Range("O13").Value = 21000
$ws.Range("P13").Value = 20500
$ws.Range("S13").Value = 1025

# Row 14
$ws.Range("D14").Value = "10/05/2021"
$ws.Range("M14").Value = 200
$ws.Range("N14").Value = 19000
$ws.Range("O14").Value = 20000
$ws.Range("P14").Value = 19500
$ws.Range("S14").Value = 975

# Row 15
$ws.Range("D15").Value = "06/03/2021"
$ws.Range("M15").Value = 160
$ws.Range("N15").Value = 19000
$ws.Range("P15").Value = 19500
$ws.Range("S15").Value = 975

# Row 16
$ws.Range("D16").Value = "09/27/2021"
$ws.Range("M16").Value = 100
$ws.Range("N16").Value = 20000
$ws.Range("O16").Value = 21000
$ws.Range("P16").Value = 20500
$ws.Range("S16").Value = 1025

# Row 17
$ws.Range("D17").Value = "05/17/2021"
$ws.Range("N17").Value = 19500
$ws.Range("O17").Value = 20000
$ws.Range("P17").Value = 19750
$ws.Range("S17").Value = 988

# Row 18
$ws.Range("D18").Value = "05/10/2021"
$ws.Range("N18").Value = 19500
$ws.Range("O18").Value = 20000
$ws.Range("P18").Value = 19750
$ws.Range("S18").Value = 988

# Row 19
$ws.Range("D19").Value = "08/12/2021"
$ws.Range("M19").Value = 160

# Row 20
$ws.Range("D20").Value = "05/27/2021"
$ws.Range("N20").Value = 19500
$ws.Range("O20").Value = 20000
$ws.Range("P20").Value = 19750
$ws.Range("S20").Value = 988

# Row 21
$ws.Range("D21").Value = "04/29/2021"
$ws.Range("N21").Value = 20000
$ws.Range("O21").Value = 21000
$ws.Range("P21").Value = 20500
$ws.Range("S21").Value = 1025

# Row 22
$ws.Range("D22").Value = "05/20/2021"
$ws.Range("M22").Value = 100
$ws.Range("N22").Value = 19500
$ws.Range("O22").Value = 20000
$ws.Range("P22").Value = 19750
$ws.Range("S22").Value = 988

# Row 24
$ws.Range("D24").Value = "06/18/2021"
$ws.Range("M24").Value = 100

# Row 25
$ws.Range("D25").Value = "10/04/2021"
$ws.Range("M25").Value = 40
$ws.Range("N25").Value = 19500
$ws.Range("O25").Value = 20000
$ws.Range("P25").Value = 19750
$ws.Range("S25").Value = 988

# Row 26
$ws.Range("D26").Value = "09/23/2021"
$ws.Range("M26").Value = 100
$ws.Range("N26").Value = 19500
$ws.Range("O26").Value = 20000
$ws.Range("P26").Value = 19750
$ws.Range("S26").Value = 988

# Row 27
$ws.Range("D27").Value = "09/06/2021"
$ws.Range("M27").Value = 160
$ws.Range("N27").Value = 20000
$ws.Range("O27").Value = 21000
$ws.Range("P27").Value = 20500
$ws.Range("S27").Value = 1025

# Row 28
$ws.Range("D28").Value = "04/15/2021"
$ws.Range("N28").Value = 18000
$ws.Range("O28").Value = 19000
$ws.Range("P28").Value = 18500
$ws.Range("S28").Value = 925

# Row 29
$ws.Range("D29").Value = "08/10/2021"
$ws.Range("N29").Value = 20000
$ws.Range("O29").Value = 21000
$ws.Range("P29").Value = 20500
$ws.Range("S29").Value = 1025

# Row 30
$ws.Range("D30").Value = "05/19/2021"
$ws.Range("M30").Value = 200
$ws.Range("N30").Value = 19000
$ws.Range("O30").Value = 20000
$ws.Range("P30").Value = 19500
$ws.Range("S30").Value = 975
